$d = $word.ActiveDocument
$r = $d.Content
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Criar estrutura de rotas</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Criar modulo de rotas (for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>child</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) e importar no </w:t></w:r><w:r><w:t>modulo principal musicas</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Criar </w:t></w:r><w:r><w:t>modulo</w:t></w:r><w:r><w:t xml:space="preserve"> de rotas principal</w:t></w:r><w:r><w:t xml:space="preserve"> (for root)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>e importar o modulo principal musicas</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Importar o modulo de rotas principal (for root) </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">no  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>app</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>.module</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Instalar o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>font</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>awesome</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para eventuais utilizações de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>icones</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Criar menu lateral contend</w:t></w:r><w:r><w:t xml:space="preserve">o um Link chamado “Repertorio” </w:t></w:r><w:r><w:t>ao qual irá abrir a listagem-musicas</w:t></w:r></w:p>

'@
$r.InsertXML($xml)
